# Refresh the live crypto snapshot (Price + Volume(1h)) pulled in the
# "Updated cryptos list" GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new display text, exactly as scraped for this refresh.
$updates = [ordered]@{
    "D2" = '27.080.25'
    "E2" = '  -2.25%  '
    "D3" = '1.865.42'
    "E3" = '  -2.04%  '
    "E4" = '  +0.03%  '
    "D5" = '306.62'
    "E5" = '  -1.88%  '
    "D6" = '1.000'
    "E6" = '  +0.08%  '
    "D7" = '0.5117'
    "E7" = '  -1.54%  '
    "D8" = '0.3749'
    "E8" = '  -0.74%  '
    "E9" = '  -1.23%  '
    "D10" = '0.8878'
    "E10" = '  -1.39%  '
    "D11" = '20.67'
    "E11" = '  -2.96%  '
    "D12" = '0.07572'
    "E12" = '  -0.74%  '
    "D13" = '1.857.82'
    "E13" = '  -2.38%  '
    "D14" = '5.306'
    "E14" = '  -2.57%  '
    "D15" = '89.49'
    "E15" = '  -2.78%  '
    "E16" = '  +0.03%  '
    "D17" = '0.000008435'
    "E17" = '  -3.00%  '
    "D18" = '14.11'
    "E18" = '  -2.55%  '
    "D19" = '0.9998'
    "E19" = '  +0.02%  '
    "D20" = '27.103.12'
    "E20" = '  -2.32%  '
    "D21" = '5.033'
    "E21" = '  -2.13%  '
    "D22" = '2.098.99'
    "E22" = '  -1.00%  '
    "E23" = '  -2.70%  '
    "D24" = '6.457'
    "E24" = '  -1.87%  '
    "D25" = '1.843'
    "E25" = '  -2.08%  '
    "D26" = '147.74'
    "E26" = '  -3.66%  '
    "D27" = '17.97'
    "E27" = '  -1.92%  '
    "D28" = '2.106'
    "E28" = '  -2.59%  '
    "D29" = '112.80'
    "E29" = '  -1.45%  '
    "D30" = '4.663'
    "E30" = '  -3.89%  '
    "D31" = '4.698'
    "E31" = '  -3.23%  '
    "D32" = '0.09096'
    "E32" = '  +1.33%  '
    "D33" = '0.05136'
    "E33" = '  -2.74%  '
    "D34" = '3.048'
    "E34" = '  -3.90%  '
    "E35" = '  -6.11%  '
    "D36" = '0.7268'
    "E36" = '  -5.65%  '
    "E37" = '  -1.98%  '
    "D38" = '2.494'
    "E38" = '  -5.59%  '
    "D39" = '3.038'
    "E39" = '  -0.82%  '
    "D40" = '1.077'
    "E40" = '  -1.63%  '
    "D41" = '0.5310'
    "E41" = '  -3.67%  '
    "D42" = '6.565'
    "E42" = '  -1.47%  '
    "D43" = '116.93'
    "E43" = '  +2.12%  '
    "D44" = '8.267'
    "E44" = '  -2.91%  '
    "D45" = '0.1469'
    "E45" = '  -2.71%  '
    "D46" = '0.9999'
    "E46" = '  +0.10%  '
    "E47" = '  -3.71%  '
    "D48" = '9.972'
    "E48" = '  -4.39%  '
    "D49" = '1.568'
    "E49" = '  -2.77%  '
    "D50" = '36.53'
    "D51" = '63.85'
    "E51" = '  -4.17%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $cell = $ws.Range($cellRef)
    if ($cellRef.StartsWith("D") -and $value -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        # Price column: force text storage so values like "1.000" or
        # "0.000008435" survive as the literal scraped string instead of
        # being normalized into a number (Excel would otherwise coerce
        # "1.000" -> 1, dropping the trailing zeros).
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
